$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = "Influenza Rapid"
$ws.Range("C2").Value = 3925638
$ws.Range("D2").Value = 32023011000653
$ws.Range("E2").Value = "Positive"
$ws.Range("F2").Value = "Categorical"

# Delete rows 3 through 9 (old additional data rows)
$ws.Range("A3:F9").EntireRow.Delete()
